$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Restyle column E (rows 3-5) to match row 2's highlighted/2-decimal style,
#     by copying the format from E2 (avoids creating any new style record).
$ws.Range("E2").Copy() | Out-Null
$ws.Range("E3:E5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Clear the (border-only) formatting from A2:A5 - these cells become
#     plain/default styled in the new data.
$ws.Range("A2:A5").ClearFormats()

# --- Write the new data set (rows 2-5), keep existing B/C/D/E styles as-is.
$ws.Range("A2").Value = 2148196
$ws.Range("B2").Value = 41466
$ws.Range("C2").Value = 7
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 3.5

$ws.Range("A3").Value = 2158422
$ws.Range("B3").Value = 28355
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 3.5

$ws.Range("A4").Value = 2172549
$ws.Range("B4").Value = 54587
$ws.Range("C4").Value = 5
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 3.5

$ws.Range("A5").Value = 2199581
$ws.Range("B5").Value = 14652
$ws.Range("C5").Value = 7
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 3.5

# --- Rows 6 & 7: the old data that used to live here is gone - clear the
#     values but keep the existing (highlighted) formatting untouched.
$ws.Range("A6:E7").ClearContents()

# --- Move the active selection, as recorded in the sheet view.
$ws.Range("C13").Select()
